# Update the "Förändrad" date column (C2:C10) from 45204 (2023-10-05)
# to 45207 (2023-10-08) as part of the automatic file update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
